$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("AB8").Value = 21
$ws.Range("AJ8").Value = 7
$ws.Range("AK8").Value = 15
$ws.Range("AR8").Value = 4.8
$ws.Range("G8").Value = 2.2
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 3.1
$ws.Range("L8").Value = 4.5
$ws.Range("Z8").Value = 9

# Row 9
$ws.Range("AA9").Value = 11
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 26
$ws.Range("AE9").Value = 5
$ws.Range("AK9").Value = 17
$ws.Range("AR9").Value = 5.8
$ws.Range("AS9").Value = 1.14
$ws.Range("G9").Value = 2.2
$ws.Range("H9").Value = 2.8
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 1.8
$ws.Range("L9").Value = 4.75
$ws.Range("U9").Value = 1.73
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 2.5
$ws.Range("X9").Value = 1.5

# Row 10
$ws.Range("AC10").Value = 34
$ws.Range("AR10").Value = 4.7
$ws.Range("AS10").Value = 1.19
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 2.55
$ws.Range("Q10").Value = 2.88
$ws.Range("R10").Value = 1.4
$ws.Range("W10").Value = 2.25
$ws.Range("X10").Value = 1.57

# Row 11
$ws.Range("AE11").Value = 5
$ws.Range("AF11").Value = 5.5
$ws.Range("AH11").Value = 81
$ws.Range("W11").Value = 2.25
$ws.Range("X11").Value = 1.57
$ws.Range("Z11").Value = 13

# Row 13
$ws.Range("AB13").Value = 17
$ws.Range("AC13").Value = 23
$ws.Range("AG13").Value = 26
$ws.Range("AH13").Value = 126
$ws.Range("AJ13").Value = 8.5
$ws.Range("AK13").Value = 21
$ws.Range("G13").Value = 1.95
$ws.Range("H13").Value = 2.8
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 2.88
$ws.Range("L13").Value = 6
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 1.67
$ws.Range("P13").Value = 2.1
$ws.Range("Q13").Value = 3.4
$ws.Range("R13").Value = 1.33
$ws.Range("U13").Value = 1.73
$ws.Range("V13").Value = 2.08
$ws.Range("Y13").Value = 4.5
$ws.Range("Z13").Value = 7.5

# Row 14
$ws.Range("O14").Value = 2
$ws.Range("P14").Value = 1.73
$ws.Range("S14").Value = 11
$ws.Range("T14").Value = 1.05

# Row 21
$ws.Range("AA21").Value = 11
$ws.Range("AB21").Value = 26
$ws.Range("G21").Value = 2.63
$ws.Range("I21").Value = 2.88
$ws.Range("M21").Value = 1.13
$ws.Range("N21").Value = 6
$ws.Range("Q21").Value = 2.88
$ws.Range("R21").Value = 1.4
$ws.Range("Y21").Value = 6

# Row 22
$ws.Range("AB22").Value = 9.5
$ws.Range("AG22").Value = 23
$ws.Range("AJ22").Value = 13
$ws.Range("AK22").Value = 34
$ws.Range("AM22").Value = 81
$ws.Range("AN22").Value = 51
$ws.Range("AP22").Value = 1.7
$ws.Range("AQ22").Value = 2.12
$ws.Range("AR22").Value = 3.45
$ws.Range("AS22").Value = 1.3
$ws.Range("G22").Value = 1.48
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 2.1
$ws.Range("L22").Value = 7.5
$ws.Range("M22").Value = 1.07
$ws.Range("N22").Value = 9
$ws.Range("W22").Value = 2.38
$ws.Range("X22").Value = 1.53
$ws.Range("Y22").Value = 5.5

# Row 45
$ws.Range("AR45").Value = 2.05
$ws.Range("AS45").Value = 1.8
$ws.Range("N45").Value = 15
$ws.Range("Q45").Value = 1.62
$ws.Range("R45").Value = 2.25

# Row 46
$ws.Range("AE46").Value = 10
$ws.Range("J46").Value = 3.25
$ws.Range("L46").Value = 3.25
$ws.Range("N46").Value = 10
$ws.Range("O46").Value = 1.3
$ws.Range("P46").Value = 3.4
$ws.Range("R46").Value = 1.8

# Row 47
$ws.Range("AB47").Value = 13.5
$ws.Range("AC47").Value = 12.5
$ws.Range("AD47").Value = 23
$ws.Range("AE47").Value = 12
$ws.Range("AF47").Value = 7.3
$ws.Range("AG47").Value = 14.5
$ws.Range("AH47").Value = 60
$ws.Range("AJ47").Value = 14
$ws.Range("AK47").Value = 27
$ws.Range("AL47").Value = 14.5
$ws.Range("AM47").Value = 75
$ws.Range("AN47").Value = 40
$ws.Range("G47").Value = 1.7
$ws.Range("H47").Value = 3.7
$ws.Range("I47").Value = 4.4
$ws.Range("J47").Value = 2.22
$ws.Range("K47").Value = 2.22
$ws.Range("L47").Value = 4.6
$ws.Range("P47").Value = 3.4
$ws.Range("X47").Value = 1.98
$ws.Range("Y47").Value = 7.8
$ws.Range("Z47").Value = 8.5

# Row 77
$ws.Range("AA77").Value = 9.25
$ws.Range("AB77").Value = 23
$ws.Range("AC77").Value = 22
$ws.Range("AD77").Value = 40
$ws.Range("AE77").Value = 5.2
$ws.Range("AF77").Value = 5.5
$ws.Range("AG77").Value = 16
$ws.Range("AH77").Value = 100
$ws.Range("AJ77").Value = 8.25
$ws.Range("AK77").Value = 18
$ws.Range("AL77").Value = 12.5
$ws.Range("AM77").Value = 55
$ws.Range("AN77").Value = 37
$ws.Range("AO77").Value = 50
$ws.Range("G77").Value = 2.25
$ws.Range("H77").Value = 2.77
$ws.Range("I77").Value = 3.55
$ws.Range("J77").Value = 2.9
$ws.Range("K77").Value = 1.9
$ws.Range("L77").Value = 4.05
$ws.Range("O77").Value = 1.47
$ws.Range("P77").Value = 2.32
$ws.Range("Q77").Value = 2.35
$ws.Range("R77").Value = 1.47
$ws.Range("S77").Value = 4
$ws.Range("T77").Value = 1.16
$ws.Range("U77").Value = 1.5
$ws.Range("V77").Value = 2.27
$ws.Range("W77").Value = 1.98
$ws.Range("X77").Value = 1.65
$ws.Range("Y77").Value = 5.8
$ws.Range("Z77").Value = 9.75
